# Append the new mods-count data point for 2026/01/28 as row 79.
#
# The "Date" column in this sheet stores ISO-like date strings ("YYYY/MM/DD")
# as literal TEXT, not as real Excel date serials. A plain
# `Range.Value = "2026/01/28"` assignment would be auto-parsed as a date by
# Excel's input parser, so instead we:
#   1. Copy the previous data row (A78:C78) down into A79:C79. This clones
#      the existing cell formatting (centered alignment) onto the new row
#      without creating any new styles.
#   2. Overwrite B79/C79 with the new Game/ModCount values directly - these
#      aren't date-like so no special handling is required.
#   3. Overwrite A79 using a leading apostrophe (the standard Excel
#      "force text" input idiom) so "2026/01/28" is stored as text rather
#      than being reinterpreted as a date, then restore the cell to the
#      sheet's normal centered style (the apostrophe entry marks the cell
#      with a quote-prefix style, so we reset it back to "Normal" and
#      reapply the same center/center alignment used throughout the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clone formatting of the last existing row onto the new row.
$ws.Range("A78:C78").Copy($ws.Range("A79:C79"))

# 2. New Game / ModCount values (plain values, no date ambiguity).
$ws.Range("B79").Value = "逃离鸭科夫"
$ws.Range("C79").Value = 1161

# 3. New Date value, forced to text, with the row's normal alignment restored.
$ws.Range("A79").Value = "'2026/01/28"
$ws.Range("A79").Style = "Normal"
$ws.Range("A79").HorizontalAlignment = -4108
$ws.Range("A79").VerticalAlignment = -4108
